# Apply the commit's changes to the presentation.
#
# 1) Slide 7's notes body placeholder: replace the "十四五规划..." boilerplate
#    with three new bullet paragraphs about the efficiency/functionality/
#    performance review notes.
# 2) Slide 8's notes body placeholder: the "十四五规划..." boilerplate is
#    removed entirely (left blank).
# 3) Slide 9's notes body placeholder: replace the "十四五规划..." boilerplate
#    with four new bullet paragraphs about offloading work to the GPU.
# 4) Slide 8's table: fix a typo, "Divide independent set" -> "Divide
#    independent sets".

$p = $ppt.ActivePresentation

# --- 1) Slide 7 notes: new "problems found" bullets -----------------------
$slide7 = $p.Slides.Item(7)
$notes7 = $slide7.NotesPage
$body7 = $notes7.Shapes.Item(2)
$body7.TextFrame.TextRange.Text = "1)确实解决了效率问题`n2）功能性缺失：信道、移动模型、协议栈、simulator实现、物理机交互`n3）性能方面：线程竞争导致伪共享问题，忽略了GPU的存在"

# --- 2) Slide 8 notes: boilerplate removed, left empty --------------------
$slide8 = $p.Slides.Item(8)
$notes8 = $slide8.NotesPage
$body8 = $notes8.Shapes.Item(2)
$body8.TextFrame.TextRange.Text = ""

# --- 3) Slide 9 notes: new "GPU offload" bullets ---------------------------
$slide9 = $p.Slides.Item(9)
$notes9 = $slide9.NotesPage
$body9 = $notes9.Shapes.Item(2)
$body9.TextFrame.TextRange.Text = "1)分离控制逻辑和计算分离，将计算任务卸载到GPU；`n2)阶段性卸载任务到GPU，然后GPU离线执行任务；`n3）kernel函数减少分支，GPU对于不执行的分支的处理是空转等待`n4）性能优化：缓存优化（解决伪共享问题）、数据交换、少用第三方库"

# --- 4) Slide 8 table: "Divide independent set" -> "...sets" --------------
$tableShape = $null
for ($i = 1; $i -le $slide8.Shapes.Count; $i++) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}
$tbl = $tableShape.Table
$cell = $tbl.Cell(1, 2)
$cell.Shape.TextFrame.TextRange.Text = "Divide independent sets"
